$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.913.40"
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("D3").Value = "2.496.48"
$ws.Range("E3").Value = "  -0.97%  "
$ws.Range("E4").Value = "  +0.04%  "
$c = $ws.Range("D5")
$c.Value = "'535.39"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.20%  "
$ws.Range("E6").Value = "  -1.96%  "
$c = $ws.Range("D7")
$c.Value = "'0.997"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("E8").Value = "  +0.64%  "
$ws.Range("D9").Value = "2.519.67"
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("E10").Value = "  +1.59%  "
$ws.Range("E11").Value = "  -0.26%  "
$c = $ws.Range("D12")
$c.Value = "'5.35"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -1.21%  "
$ws.Range("E13").Value = "  -2.37%  "
$ws.Range("D14").Value = "2.947.98"
$ws.Range("E14").Value = "  -0.52%  "
$c = $ws.Range("D15")
$c.Value = "'23.24"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +1.15%  "
$ws.Range("D16").Value = "58.874.38"
$ws.Range("E16").Value = "  -0.43%  "
$ws.Range("E17").Value = "  -0.39%  "
$ws.Range("D18").Value = "2.517.76"
$ws.Range("E18").Value = "  -0.04%  "
$c = $ws.Range("D19")
$c.Value = "'11.05"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +1.28%  "
$ws.Range("E20").Value = "  +0.81%  "
$c = $ws.Range("D21")
$c.Value = "'325.10"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +1.05%  "
$ws.Range("E22").Value = "  +0.44%  "
$ws.Range("E23").Value = "  +0.60%  "
$c = $ws.Range("D24")
$c.Value = "'64.83"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +3.15%  "
$ws.Range("E25").Value = "  -1.13%  "
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("E27").Value = "  +1.07%  "
$c = $ws.Range("D28")
$c.Value = "'7.61"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -1.95%  "
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("E30").Value = "  +0.84%  "
$ws.Range("E31").Value = "  -1.42%  "
$c = $ws.Range("D32")
$c.Value = "'167.70"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +4.44%  "
$ws.Range("E33").Value = "  +4.00%  "
$c = $ws.Range("D34")
$c.Value = "'0.998"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.13%  "
$c = $ws.Range("D35")
$c.Value = "'1.41"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -3.50%  "
$ws.Range("E36").Value = "  +0.52%  "
$c = $ws.Range("D37")
$c.Value = "'4.13"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -1.77%  "
$c = $ws.Range("D38")
$c.Value = "'1.57"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -0.97%  "
$c = $ws.Range("D39")
$c.Value = "'36.84"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -0.29%  "
$c = $ws.Range("D40")
$c.Value = "'0.832"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +3.35%  "
$c = $ws.Range("D41")
$c.Value = "'3.64"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +0.08%  "
$c = $ws.Range("D42")
$c.Value = "'5.30"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +0.59%  "
$c = $ws.Range("D43")
$c.Value = "'280.85"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -0.91%  "
$ws.Range("E44").Value = "  -0.29%  "
$ws.Range("E45").Value = "  +1.70%  "
$c = $ws.Range("D46")
$c.Value = "'10.88"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +0.21%  "
$c = $ws.Range("D47")
$c.Value = "'127.81"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +4.07%  "
$ws.Range("E48").Value = "  -0.06%  "
$c = $ws.Range("D49")
$c.Value = "'0.0512"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +0.61%  "
$ws.Range("E50").Value = "  -0.02%  "
$c = $ws.Range("D51")
$c.Value = "'17.36"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -0.29%  "
